# ---------------------------------------------------------------------------
# queryQy.py pytest-class demo sheet: add a new "JD查询企业信息" worksheet with
# sample request/response fixture data, mirroring the existing
# mtPossession/juhe/register sheets' layout.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- add the new worksheet as the last tab -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "JD查询企业信息"

# --- column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.375
$ws.Columns.Item(2).ColumnWidth = 20.5
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 17.5

# --- header row (row 1) --------------------------------------------------
$headerFill = 5296274          # RGB(146,208,80) == fill used elsewhere in the workbook

$hdrA = $ws.Range("A1")
$hdrA.NumberFormat = "@"
$hdrA.Font.Name = "等线"
$hdrA.Font.Size = 11
$hdrA.Interior.Color = $headerFill
$hdrA.Borders.LineStyle = 1
$hdrA.Borders.Weight = 2
$hdrA.VerticalAlignment = -4108
$hdrA.Value = "test_name"

$hdrRestGeneral = $ws.Range("B1,D1:F1")
$hdrRestGeneral.Font.Name = "等线"
$hdrRestGeneral.Font.Size = 11
$hdrRestGeneral.Interior.Color = $headerFill
$hdrRestGeneral.Borders.LineStyle = 1
$hdrRestGeneral.Borders.Weight = 2
$hdrRestGeneral.VerticalAlignment = -4108

$hdrC = $ws.Range("C1")
$hdrC.NumberFormat = "@"
$hdrC.Font.Name = "等线"
$hdrC.Font.Size = 11
$hdrC.Interior.Color = $headerFill
$hdrC.Borders.LineStyle = 1
$hdrC.Borders.Weight = 2
$hdrC.VerticalAlignment = -4108

$ws.Range("B1").Value = "test_description"
$ws.Range("C1").Value = "enterpriseName"
$ws.Range("D1").Value = "ssoId"
$ws.Range("E1").Value = "platform"
$ws.Range("F1").Value = "code"

$ws.Rows.Item(1).RowHeight = 12.75

# --- data rows (rows 2-7) : A/B/D/F use the text number format ----------
$dataText = $ws.Range("A2:B7,D2:D7,F2:F7")
$dataText.NumberFormat = "@"

$enterpriseName = "长沙市天心区祝博士教育信息咨询服务部"

$data = @(
    @("001", "查询成功",   $enterpriseName, "1445898791152850915", "qyd", "10000"),
    @("002", "参数错误",   $enterpriseName, "1445898791152850915", "qyd", "50003"),
    @("003", "数据库操作错误", $enterpriseName, "1445898791152850915", "qyd", "50004"),
    @("004", "查询服务已关闭", $enterpriseName, "1445898791152850915", "qyd", "50006"),
    @("005", "调用第三方其他错误", $enterpriseName, "1445898791152850915", "qyd", "50007"),
    @("006", "JD业务失败", $enterpriseName, "1445898791152850915", "qyd", "50008")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# odd data rows (3,5,7) render the enterprise name in a muted grey Arial font
$ws.Range("C3,C5,C7").Font.Name = "Arial"
$ws.Range("C3,C5,C7").Font.Size = 11
$ws.Range("C3,C5,C7").Font.Color = 3355443   # RGB(51,51,51) == FF333333

# --- trailing blank rows (8-17), mirroring the author's placeholder rows --
$ws.Range("A8:D8,F8:F8").NumberFormat = "@"
$ws.Range("A9:D17,F9:F17").NumberFormat = "@"

# --- sheet view: fresh selection on the new sheet, matches author's edit --
$ws.Range("H10").Select()

# register's sheetView previously held the "active tab" marker; after adding
# the new sheet it moves with the active tab automatically.
$ws.Activate()
